# Apply "replacing files with master branch files" changes to the Sanity
# workbook's single sheet (TestCaseMaster), row 14 (the "All APIs" test
# case), plus the associated view-state tweak.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 currently reads:
#   A14=13  B14=All APIs  C14=Verifying Auth Login API
#   D14=com.darwinbox.mobile.allAPIs.AllApiStatusCheckTest
#   E14=allAPI_StatusCheck/allAPI_StatusCheck_TestData.xlsx
#   F14=test2_instance
#   G14=1-5,7,9,11-12,14-15,18,20-23,27-40,42,44-48,50,51
#
# C14/D14/E14 keep their existing text (only the shared-string table's
# internal ordering shifts upstream, which is not user-visible); the real
# content edits are F14, G14, and the brand-new H14:
#   F14=wild1_instance
#   G14=1
#   H14=1-10,12-16,20-37,39-43,46-54,56-60,62-63   (new column)

$ws.Range("F14").Value = "wild1_instance"

# G14 and the new H14 hold digit-lead strings ("1", and a range list) that
# must stay text (quote-prefixed), matching how the rest of the sheet
# stores its numeric-looking TestDataRow/TCID values.
$ws.Range("G14").Value = "'1"
$ws.Range("H14").Value = "'1-10,12-16,20-37,39-43,46-54,56-60,62-63"

# Keep the selection where it was (G14 in the source, still valid there).
$ws.Range("G14").Select() | Out-Null

# The source view also scrolls the pane so column C becomes the leftmost
# visible column (was column B); best-effort, mirrored via the window's
# scroll position.
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 3
$aw.ScrollRow = 1
